# Fix off-by-one "index" numbering in a few sheets (per commit: "fix some dirty files").
# Each affected row has its row-index value (column A) and the duplicated "index"
# column (the last column on the sheet) decremented by 1.

$wb = $excel.ActiveWorkbook

# --- 股票 (stock) sheet: rows 2-3, index column = N ---
$ws = $wb.Worksheets.Item("股票")
$ws.Range("A2").Value = 69
$ws.Range("N2").Value = 69
$ws.Range("A3").Value = 70
$ws.Range("N3").Value = 70

# --- 具有相當價值之財產 (valuable property) sheet: rows 2-6, index column = L ---
$ws = $wb.Worksheets.Item("具有相當價值之財產")
$ws.Range("A2").Value = 90
$ws.Range("L2").Value = 90
$ws.Range("A3").Value = 91
$ws.Range("L3").Value = 91
$ws.Range("A4").Value = 92
$ws.Range("L4").Value = 92
$ws.Range("A5").Value = 94
$ws.Range("L5").Value = 94
$ws.Range("A6").Value = 95
$ws.Range("L6").Value = 95

# --- 保險 (insurance) sheet: rows 2-9, index column = K ---
$ws = $wb.Worksheets.Item("保險")
$ws.Range("A2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("A3").Value = 101
$ws.Range("K3").Value = 101
$ws.Range("A4").Value = 102
$ws.Range("K4").Value = 102
$ws.Range("A5").Value = 104
$ws.Range("K5").Value = 104
$ws.Range("A6").Value = 105
$ws.Range("K6").Value = 105
$ws.Range("A7").Value = 106
$ws.Range("K7").Value = 106
$ws.Range("A8").Value = 107
$ws.Range("K8").Value = 107
$ws.Range("A9").Value = 108
$ws.Range("K9").Value = 108
